$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.734.81"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "2.079.22"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  +0.45%  "
$__style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = $__style
$ws.Range("E7").Value = "  -0.09%  "
$__style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.07"
$ws.Range("D8").Style = $__style
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("D12").Value = "2.385.76"
$ws.Range("E12").Value = "  +0.87%  "
$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.49"
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = "  +0.59%  "
$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.96"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("E15").Value = "  -1.57%  "
$__style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.26"
$ws.Range("D16").Style = $__style
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "2.083.38"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "37.687.69"
$ws.Range("E18").Value = "  +1.22%  "
$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.18"
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = "  -2.89%  "
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.73"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("E21").Value = "  +0.95%  "
$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.81"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("E25").Value = "  -0.62%  "
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.96"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = "  +2.27%  "
$ws.Range("E27").Value = "  +11.31%  "
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("E29").Value = "  +0.31%  "
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.39"
$ws.Range("D30").Style = $__style
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  +0.86%  "
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.65"
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = "  +3.87%  "
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("E36").Value = "  +2.44%  "
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.38"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = "  +5.41%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  -1.88%  "
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0990"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = "  +6.81%  "
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.96"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = "  +3.17%  "
$ws.Range("E42").Value = "  -0.67%  "
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.36"
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").Value = "1.456.60"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("E47").Value = "  +3.14%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.44"
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = "  +3.75%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.60"
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("E50").Value = "  +2.04%  "
$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.59"
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = "  +7.55%  "
